$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert the two new rows at their final target positions.
#    Row 5  -> "Vista 4"          (pushes old rows 5..12 down to 6..13)
#    Row 10 -> "Procedimiento 5"  (pushes old rows 10..13 down to 11..14)
# ---------------------------------------------------------------------------
$ws.Rows.Item(5).Insert()
$ws.Rows.Item(10).Insert()

# New rows created by Insert() copy the formatting of the row above them,
# but not the custom row height - restore it explicitly.
$ws.Rows.Item(5).RowHeight = 30
$ws.Rows.Item(10).RowHeight = 30

# ---------------------------------------------------------------------------
# 2. Fill in the content for the brand new rows.
# ---------------------------------------------------------------------------
# Row 5 - Vista 4
$ws.Range("A5").Value = "Vista 4"
$ws.Range("B5").Value = "Mostrar articulos con stock bajo"
$ws.Range("D5").Value = "OK"
$ws.Range("E5").Value = "Leandro / Facundo"

# Row 10 - Procedimiento 5
$ws.Range("A10").Value = "Procedimiento 5"
$ws.Range("B10").Value = "Realizar alta de Marca"
$ws.Range("C10").Value = "Realizar insert de marca"
$ws.Range("D10").Value = "OK"
$ws.Range("E10").Value = "Leandro Correa"

# ---------------------------------------------------------------------------
# 3. Update the existing "Trigger 1" row (now row 11) with its new data.
# ---------------------------------------------------------------------------
$ws.Range("B11").Value = "Notificación de productos con stok bajo"
$ws.Range("D11").Value = "OK"
$ws.Range("E11").Value = "Facundo"

# ---------------------------------------------------------------------------
# 4. Update the "Video Demostrativo" row (now row 13): Cumplido -> EN CURSO
# ---------------------------------------------------------------------------
$ws.Range("D13").Value = "EN CURSO"

# ---------------------------------------------------------------------------
# 5. Remove the trailing blank row (old row 12, now shifted to row 14).
# ---------------------------------------------------------------------------
$ws.Rows.Item(14).Delete()

# ---------------------------------------------------------------------------
# 6. Column C no longer uses wrap-text formatting anywhere.
# ---------------------------------------------------------------------------
$ws.Range("C2:C13").WrapText = $false

# ---------------------------------------------------------------------------
# 7. Resize the table to include the new rows.
# ---------------------------------------------------------------------------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:E13"))

# ---------------------------------------------------------------------------
# 8. Column width adjustments for A and B.
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 18.02
$ws.Columns.Item(2).ColumnWidth = 67.88

# ---------------------------------------------------------------------------
# 9. Update the active cell selection to match the edited workbook.
# ---------------------------------------------------------------------------
$ws.Range("D13").Select()
